# The document originally contains a handful of blank paragraphs between
# "This is the introduction" and the "The organization of this article..."
# paragraph. The 3rd blank paragraph (overall paragraph #4) is turned into
# a short "Stylised facts" write-up: a bold heading, a blank bold line, a
# bold "Ideas:" line, and a bulleted (plain) list of table/figure ideas,
# all tagged as English (UK) instead of the surrounding English (US) text.
$d = $word.ActiveDocument

$target = $d.Paragraphs.Item(4)
$targetRange = $target.Range

$newContent = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-GB"/></w:rPr><w:t>Stylised facts</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-GB"/></w:rPr><w:t>Ideas:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Total number of SS, NS and NN PTAs</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Share of SS, NS and NN PTAs</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Total exports by S and N countries</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Share of total exports by S and N countries</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Total exports of manufactured by S and N countries</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Share of total exports of manufactured by S and N countries</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Number of products exported </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>by S and N countries</w:t></w:r></w:p>
'@

# InsertXML replaces the contents of the addressed range (here, the
# whole empty paragraph, mark included) with the supplied WordprocessingML,
# so this single call both retypes the first paragraph and appends the
# nine paragraphs that follow it.
[void]$targetRange.InsertXML($newContent)

Write-Host "Paragraph count after edit:" $d.Paragraphs.Count
